$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6309176683425903
$ws.Range("B1").Value = 2.343778371810913
$ws.Range("C1").Value = 3.19474196434021
$ws.Range("D1").Value = 3.788445472717285
$ws.Range("E1").Value = 1.330827951431274
